$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows before row 179 (existing rows 179.. shift down to 182..)
$ws.Rows("179:181").Insert()

# --- New row 179 ---
$ws.Range("A179").Value = 6
$ws.Range("B179").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C179").Value = "Metropolitana"
$ws.Range("D179").Value = 45258
$ws.Range("E179").Value = 13
$ws.Range("F179").Value = 300000000
$ws.Range("G179").Value = "Espárragos"
$ws.Range("H179").Value = "Sin especificar"
$ws.Range("I179").Value = "Banquete"
$ws.Range("J179").Value = 350
$ws.Range("K179").Value = 1600
$ws.Range("L179").Value = 1800
$ws.Range("M179").Value = 1714
$ws.Range("N179").Value = "`$/kilo"
$ws.Range("O179").Value = "Provincia de Linares"
$ws.Range("P179").Value = 1714
$ws.Range("Q179").Value = 1
$ws.Range("R179").Value = "Hortaliza"

# --- New row 180 ---
$ws.Range("A180").Value = 6
$ws.Range("B180").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C180").Value = "Metropolitana"
$ws.Range("D180").Value = 45258
$ws.Range("E180").Value = 13
$ws.Range("F180").Value = 300000000
$ws.Range("G180").Value = "Espárragos"
$ws.Range("H180").Value = "Sin especificar"
$ws.Range("I180").Value = "Primera"
$ws.Range("J180").Value = 1100
$ws.Range("K180").Value = 1200
$ws.Range("L180").Value = 1400
$ws.Range("M180").Value = 1255
$ws.Range("N180").Value = "`$/kilo"
$ws.Range("O180").Value = "Provincia de Linares"
$ws.Range("P180").Value = 1255
$ws.Range("Q180").Value = 1
$ws.Range("R180").Value = "Hortaliza"

# --- New row 181 ---
$ws.Range("A181").Value = 6
$ws.Range("B181").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C181").Value = "Metropolitana"
$ws.Range("D181").Value = 45258
$ws.Range("E181").Value = 13
$ws.Range("F181").Value = 300000000
$ws.Range("G181").Value = "Espárragos"
$ws.Range("H181").Value = "Sin especificar"
$ws.Range("I181").Value = "Segunda"
$ws.Range("J181").Value = 1100
$ws.Range("K181").Value = 1000
$ws.Range("L181").Value = 1200
$ws.Range("M181").Value = 1073
$ws.Range("N181").Value = "`$/kilo"
$ws.Range("O181").Value = "Provincia de Linares"
$ws.Range("P181").Value = 1073
$ws.Range("Q181").Value = 1
$ws.Range("R181").Value = "Hortaliza"
